{"js": "const body = context.document.body;\nconst results = body.search(\"cp -r /scratch/\", {matchCase: true});\nresults.load(\"items\");\nawait context.sync();\nconst r = results.items[0];\nconst ooxml = r.getOoxml();\nawait context.sync();\nreturn ooxml.value;\n", "ps1": "$d = $word.ActiveDocument\n$existing = $d.Bookmarks(\"_Hlk528097548\")\n$startPos = $existing.Start\n$endPos = $existing.End\n$existing.Delete()\n\n$rGoBack = $d.Range($startPos, $startPos)\n$d.Bookmarks.Add(\"_GoBack\", $rGoBack)\n\n$rHlk = $d.Range($startPos, $endPos)\n$d.Bookmarks.Add(\"_Hlk528097548\", $rHlk)\n"}
